$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "op is “=” & got input / Calculate with last operation" textbox
# (shape id 57, "TextBox 56") - it sits right after the "No input"/"A = result"
# textbox and before the curved connector that starts at the "Got input" oval.
$shp = $s.Shapes.Item("TextBox 56")

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Replace the two-paragraph text with just the first run's new wording.
# Setting .Text collapses the textbox down to a single paragraph (dropping the
# old second paragraph "Calculate with last operation "); the autosized
# height recalculates for us because the shape has <a:spAutoFit/>.
$quote1 = [char]0x201C
$quote2 = [char]0x201D
$tr.Text = "op is " + $quote1 + "=" + $quote2 + " & got "

# Add "no input" as its own trailing run (matches the authored edit).
[void]$tr.InsertAfter("no input")

# The textbox uses wrap="none" with spAutoFit, so PowerPoint also shrinks the
# width to fit the now-shorter text (height already auto-adjusted above).
$shp.Width = 188.85142517089844
